$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths: merge A:C into a single uniform width ---
$ws.Range("A1:C1").ColumnWidth = 32.63

# --- New column S: year 2022 header + its index value, inheriting the
#     formatting already used by the neighbouring column R ---
$ws.Columns("S:S").Insert()
$ws.Range("S3").Value = 2022
$ws.Range("S4").Value = 0.071025550219041236

# --- Update the saved selection/active cell ---
$ws.Range("F14").Select()
